$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Date of Birth"

# Data row
$ws.Range("A2").Value = "Ruhullah"
$ws.Range("B2").Value = 36912
$ws.Range("B2").NumberFormat = "mm-dd-yy"

# Widen the date column to fit its contents (matches the author's "bestFit" column)
$ws.Columns.Item(2).ColumnWidth = 9.5

# Leave the selection where the author left it after entering the data
$ws.Range("B3").Select() | Out-Null
